$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.179.80"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "3.428.13"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'406.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").Value = "'132.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.57%  "

$ws.Range("D7").Value = "'0.595"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.692"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.95%  "

$ws.Range("E10").Value = "  +6.75%  "

$ws.Range("D11").Value = "'42.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").Value = "'8.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("D15").Value = "3.412.41"
$ws.Range("E15").Value = "  +1.86%  "

$ws.Range("D16").Value = "'11.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").Value = "62.174.08"
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").Value = "'0.0000146"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.59%  "

$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").Value = "'84.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.80%  "

$ws.Range("D22").Value = "'311.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("E23").Value = "  -2.98%  "

$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "'29.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").Value = "'8.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.25%  "

$ws.Range("D28").Value = "'7.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.90%  "

$ws.Range("E29").Value = "  +5.99%  "

$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("D31").Value = "'43.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("D33").Value = "'11.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.61%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "'51.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.15%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("D39").Value = "'3.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.71%  "

$ws.Range("E40").Value = "  +11.94%  "

$ws.Range("D41").Value = "'143.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.25%  "

$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("E45").Value = "  -0.78%  "

$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "'21.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.80%  "

$ws.Range("D48").Value = "2.104.49"
$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("E49").Value = "  -1.94%  "

$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("E51").Value = "  +19.06%  "
